$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "66.925.49"
$ws.Range("E2").Value = "  -3.67%  "

Set-TextValue $ws "D3" "3.532.76"
$ws.Range("E3").Value = "  -3.86%  "

$ws.Range("E4").Value = "  +0.02%  "

Set-TextValue $ws "D5" "606.87"
$ws.Range("E5").Value = "  -5.24%  "

Set-TextValue $ws "D6" "154.04"
$ws.Range("E6").Value = "  -3.39%  "

Set-TextValue $ws "D7" "3.530.39"
$ws.Range("E7").Value = "  -3.86%  "

$ws.Range("E8").Value = "  +0.07%  "

Set-TextValue $ws "D9" "0.485"
$ws.Range("E9").Value = "  -2.46%  "

Set-TextValue $ws "D10" "0.141"
$ws.Range("E10").Value = "  -2.48%  "

Set-TextValue $ws "D11" "6.83"
$ws.Range("E11").Value = "  -3.87%  "

Set-TextValue $ws "D12" "0.430"
$ws.Range("E12").Value = "  -3.70%  "

$ws.Range("E13").Value = "  -4.44%  "

Set-TextValue $ws "D14" "4.132.57"
$ws.Range("E14").Value = "  -3.81%  "

Set-TextValue $ws "D15" "31.99"
$ws.Range("E15").Value = "  -2.10%  "

Set-TextValue $ws "D16" "3.518.17"
$ws.Range("E16").Value = "  -4.07%  "

Set-TextValue $ws "D17" "66.954.60"
$ws.Range("E17").Value = "  -3.61%  "

$ws.Range("E18").Value = "  +0.78%  "

Set-TextValue $ws "D19" "6.37"
$ws.Range("E19").Value = "  -1.91%  "

Set-TextValue $ws "D20" "15.45"
$ws.Range("E20").Value = "  -3.22%  "

Set-TextValue $ws "D21" "450.04"
$ws.Range("E21").Value = "  -3.57%  "

Set-TextValue $ws "D22" "9.33"
$ws.Range("E22").Value = "  -5.59%  "

Set-TextValue $ws "D23" "0.638"
$ws.Range("E23").Value = "  -1.60%  "

Set-TextValue $ws "D24" "79.03"
$ws.Range("E24").Value = "  -0.28%  "

Set-TextValue $ws "D25" "3.674.47"
$ws.Range("E25").Value = "  -3.85%  "

$ws.Range("E26").Value = "  -0.02%  "

Set-TextValue $ws "D27" "0.0000123"
$ws.Range("E27").Value = "  -2.23%  "

Set-TextValue $ws "D28" "10.26"
$ws.Range("E28").Value = "  -5.76%  "

Set-TextValue $ws "D29" "8.32"
$ws.Range("E29").Value = "  -7.93%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws "D30" "2.55"
$ws.Range("E30").Value = "  -2.97%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws "D31" "1.68"
$ws.Range("E31").Value = "  -1.11%  "

Set-TextValue $ws "D32" "0.999"
$ws.Range("E32").Value = "  +0.01%  "

Set-TextValue $ws "D33" "25.91"
$ws.Range("E33").Value = "  -3.48%  "

Set-TextValue $ws "D34" "1.89"
$ws.Range("E34").Value = "  -5.16%  "

$ws.Range("E35").Value = "  -3.99%  "

Set-TextValue $ws "D36" "6.18"
$ws.Range("E36").Value = "  -4.24%  "

Set-TextValue $ws "D37" "3.529.14"

Set-TextValue $ws "D38" "8.11"
$ws.Range("E38").Value = "  -4.23%  "

$ws.Range("E39").Value = "  +0.05%  "

Set-TextValue $ws "D40" "1.00"
$ws.Range("E40").Value = "  -0.02%  "

Set-TextValue $ws "D41" "176.08"
$ws.Range("E41").Value = "  -0.57%  "

Set-TextValue $ws "D42" "5.61"
$ws.Range("E42").Value = "  -4.66%  "

$ws.Range("E43").Value = "  -3.49%  "

Set-TextValue $ws "D44" "0.0874"
$ws.Range("E44").Value = "  -2.84%  "

Set-TextValue $ws "D45" "0.891"
$ws.Range("E45").Value = "  -3.67%  "

Set-TextValue $ws "D46" "45.79"
$ws.Range("E46").Value = "  -2.24%  "

Set-TextValue $ws "D47" "28.24"
$ws.Range("E47").Value = "  +2.99%  "

$ws.Range("E48").Value = "  -1.52%  "

Set-TextValue $ws "D49" "1.23"
$ws.Range("E49").Value = "  -1.05%  "

Set-TextValue $ws "D50" "7.65"
$ws.Range("E50").Value = "  -2.48%  "

Set-TextValue $ws "D51" "1.03"
$ws.Range("E51").Value = "  -3.46%  "
